$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add a new day column H ("21-jun") with 24 hourly prices
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the header formatting from G1 (bold / bordered / centered style) onto
# H1, then set its text so the new header matches the look of the others.
$wsPrix.Range("G1").Copy()
$wsPrix.Range("H1").PasteSpecial(-4122)
$wsPrix.Range("H1").Value = "21-jun"

$hValues = @(
    118.53,
    110.48,
    108.15,
    98.03,
    81.62,
    92.27,
    98.79000000000001,
    90.73999999999999,
    75.93000000000001,
    23.4,
    2.5,
    0,
    0,
    0,
    0,
    0,
    0.1,
    60.4,
    98.64,
    127.53,
    141.99,
    134.93,
    138.81,
    119.48
)

for ($i = 0; $i -lt $hValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 8).Value = $hValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 5 with the 2025-06-19 closing price
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A5").Value = "'2025-06-19"
$wsGaz.Range("A4").Copy()
$wsGaz.Range("A5").PasteSpecial(-4122)
$wsGaz.Range("B5").Value = 40.425

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 5 with the 2025-06-19 closing price
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A5").Value = "'2025-06-19"
$wsCo2.Range("A4").Copy()
$wsCo2.Range("A5").PasteSpecial(-4122)
$wsCo2.Range("B5").Value = 71
